# Split backtesting in multiple steps common to all strategies.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: From/To dates move to Nov, TP/SL % updated
$ws.Range("D2").Value = 44501
$ws.Range("H2").Value = 1.2
$ws.Range("I2").Value = 1

# Row 3: From/To dates, interval, TP/SL %, strategy
$ws.Range("D3").Value = 44501
$ws.Range("F3").Value = "30m"
$ws.Range("H3").Value = 1.2
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = "MACD"

# Row 4: brand new third test case (was blank except D4/E4)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ByBit"
$ws.Range("C4").Value = "BTCUSDT"
$ws.Range("D4").Value = 44501
$ws.Range("E4").Value = 44561
$ws.Range("F4").Value = "30m"
$ws.Range("G4").Value = 10000
$ws.Range("H4").Value = 1.2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = "EarlyMACD"

# Move selection/active cell to J11 (where the user left off editing)
$ws.Range("J11").Select()
